$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2025-10-13 14:38:01", "what is the total number of orders in july", "8926 orders in July.", "5b5b4eaa-0cb2-4eab-a958-dac0b2da8b2e"),
    @("2025-10-13 15:14:06", "total number of canceled orders in july", "Total number of canceled orders in July: 257.", "0c47e2e6-2890-43fb-be0c-dfaf620e22fe"),
    @("2025-10-13 15:14:19", "in august", "Sorry, we couldn't find the data you are looking for.", "0c47e2e6-2890-43fb-be0c-dfaf620e22fe"),
    @("2025-10-13 15:20:06", "total number of canceled orders in july", "Sorry, we couldn't find the data you are looking for.", "0c47e2e6-2890-43fb-be0c-dfaf620e22fe"),
    @("2025-10-13 15:22:01", "total number of canceled orders in july", "Total number of canceled orders in July: 257.", "0c47e2e6-2890-43fb-be0c-dfaf620e22fe")
)

$startRow = 7
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
}
